# Ajout d'une vingtaine de photos pour du personnel deja dans le repertoire.
# Renseigne la colonne "photo" (colonne C) du tableau du personnel pour les
# employes qui n'avaient pas encore de photo associee, et corrige la casse
# du nom de fichier existant pour France Nadeau.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$photos = @{
    8   = "eve-baribeau-marchand.jpg"
    12  = "audrey-begin-poisson.jpg"
    17  = "marie-france-bernier.jpg"
    25  = "guylaine-brazeau.jpg"
    27  = "fanny-c-brochu.jpg"
    28  = "delphine-cado.jpg"
    29  = "siv-kham-chao.jpg"
    31  = "noemie-charest-bourbon.jpg"
    61  = "isabelle-giguere.jpg"
    88  = "benoit-mayrand.jpg"
    93  = "mathieu-murray-samuel.jpg"
    94  = "france-nadeau.jpg"
    111 = "chloe-sinotte.jpg"
    112 = "ann-mary-sotomayor.jpg"
    116 = "karine-theriault-dube.jpg"
    120 = "juliette-tirard-collet.jpg"
    121 = "frederic-tremblay.jpg"
    126 = "yan-vallee.jpg"
}

foreach ($row in $photos.Keys) {
    $ws.Cells.Item($row, 3).Value = $photos[$row]
}
